$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("line_imp")
$ws.Activate()

$ws.Range("E3").Formula = "=0.05/2"

$ws.Range("E8").Select()
